# Update TPM-derived values in the LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 20.598495
    "H2"  = 61.795485
    "I2"  = 0.187290770808477
    "J2"  = 0.187290770808477
    "K2"  = 2
    "L2"  = 0.6666666666666666
    "M2"  = 0.04182166666666667
    "N2"  = 0.125465
    "O2"  = 0.03574612430984324
    "P2"  = 0.03574612430984324
    "Q2"  = 0.861463391725
    "R2"  = 7.753170525524999
    "S2"  = 0.006694919175406179
    "T2"  = 0.006694919175406179

    "G3"  = 20.598495
    "H3"  = 61.795485
    "I3"  = 0.187290770808477
    "J3"  = 0.187290770808477
    "O3"  = 0.4725152433508619
    "P3"  = 0.4725152433508619
    "Q3"  = 11.387376730705
    "R3"  = 102.486390576345
    "S3"  = 0.08849774414593803
    "T3"  = 0.08849774414593803

    "G4"  = 20.598495
    "H4"  = 61.795485
    "I4"  = 0.187290770808477
    "J4"  = 0.187290770808477
    "M4"  = 0.5753163333333333
    "N4"  = 1.725949
    "O4"  = 0.4917386323392949
    "P4"  = 0.4917386323392949
    "Q4"  = 11.850650615585
    "R4"  = 106.655855540265
    "S4"  = 0.09209810748713283
    "T4"  = 0.09209810748713283

    "G5"  = 60.20577233333334
    "I5"  = 0.5474179306512287
    "J5"  = 0.5474179306512288
    "K5"  = 2
    "L5"  = 0.6666666666666666
    "M5"  = 0.04182166666666667
    "N5"  = 0.125465
    "O5"  = 0.03574612430984324
    "P5"  = 0.03574612430984324
    "Q5"  = 2.517905741933889
    "R5"  = 22.661151677405
    "S5"  = 0.01956806939849597
    "T5"  = 0.01956806939849597

    "G6"  = 60.20577233333334
    "I6"  = 0.5474179306512287
    "J6"  = 0.5474179306512288
    "O6"  = 0.4725152433508619
    "P6"  = 0.4725152433508619
    "Q6"  = 33.28329622735656
    "S6"  = 0.2586633167162906
    "T6"  = 0.2586633167162906

    "G7"  = 60.20577233333334
    "I7"  = 0.5474179306512287
    "J7"  = 0.5474179306512288
    "M7"  = 0.5753163333333333
    "N7"  = 1.725949
    "O7"  = 0.4917386323392949
    "P7"  = 0.4917386323392949
    "Q7"  = 34.63736418431478
    "R7"  = 311.736277658833
    "S7"  = 0.2691865445364421
    "T7"  = 0.2691865445364422

    "G8"  = 29.17709966666666
    "H8"  = 87.53129899999999
    "I8"  = 0.2652912985402942
    "J8"  = 0.2652912985402942
    "K8"  = 2
    "L8"  = 0.6666666666666666
    "M8"  = 0.04182166666666667
    "N8"  = 0.125465
    "O8"  = 0.03574612430984324
    "P8"  = 0.03574612430984324
    "Q8"  = 1.220234936559444
    "R8"  = 10.982114429035
    "S8"  = 0.009483135735941091
    "T8"  = 0.009483135735941091

    "G9"  = 29.17709966666666
    "H9"  = 87.53129899999999
    "I9"  = 0.2652912985402942
    "J9"  = 0.2652912985402942
    "O9"  = 0.4725152433508619
    "P9"  = 0.4725152433508619
    "Q9"  = 16.12984957462477
    "R9"  = 145.168646171623
    "S9"  = 0.1253541824886333
    "T9"  = 0.1253541824886333

    "G10" = 29.17709966666666
    "H10" = 87.53129899999999
    "I10" = 0.2652912985402942
    "J10" = 0.2652912985402942
    "M10" = 0.5753163333333333
    "N10" = 1.725949
    "O10" = 0.4917386323392949
    "P10" = 0.4917386323392949
    "Q10" = 16.78606199752789
    "R10" = 151.074557977751
    "S10" = 0.1304539803157198
    "T10" = 0.1304539803157198
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
